$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Edi09_2_bg")
$ws2.Range("C10").Value2 = 0.442
$ws2.Range("D10").Value2 = 0.444
$ws2.Range("E10").Value2 = 0.444
$ws2.Range("F10").Value2 = 0.445
$ws2.Range("G10").Value2 = 0.442
$ws2.Range("J10").Value2 = 0.442
$ws2.Range("K10").Value2 = 0.444
$ws2.Range("L10").Value2 = 0.444
$ws2.Range("M10").Value2 = 0.449
$ws2.Range("N10").Value2 = 0.443
$ws2.Range("O10").Value2 = 0.444
$ws2.Range("Q10").Value2 = 0.442
$ws2.Range("R10").Value2 = 0.449
$ws2.Range("B11").Value2 = 0.066
$ws2.Range("C11").Value2 = 0.106
$ws2.Range("D11").Value2 = 0.08699999999999999
$ws2.Range("E11").Value2 = 0.089
$ws2.Range("F11").Value2 = 0.082
$ws2.Range("G11").Value2 = 0.109
$ws2.Range("H11").Value2 = 0.08500000000000001
$ws2.Range("I11").Value2 = 0.092
$ws2.Range("J11").Value2 = 0.109
$ws2.Range("K11").Value2 = 0.08599999999999999
$ws2.Range("L11").Value2 = 0.08699999999999999
$ws2.Range("M11").Value2 = 0.045
$ws2.Range("N11").Value2 = 0.099
$ws2.Range("O11").Value2 = 0.08799999999999999
$ws2.Range("P11").Value2 = 0.017
$ws2.Range("Q11").Value2 = 0.045
$ws2.Range("R11").Value2 = 0.109
$ws2.Range("B16").Value2 = 46.459
$ws2.Range("C16").Value2 = 46.424
$ws2.Range("D16").Value2 = 46.44
$ws2.Range("E16").Value2 = 46.439
$ws2.Range("F16").Value2 = 46.445
$ws2.Range("G16").Value2 = 46.422
$ws2.Range("H16").Value2 = 46.443
$ws2.Range("I16").Value2 = 46.437
$ws2.Range("J16").Value2 = 46.421
$ws2.Range("K16").Value2 = 46.442
$ws2.Range("L16").Value2 = 46.441
$ws2.Range("M16").Value2 = 46.478
$ws2.Range("N16").Value2 = 46.43
$ws2.Range("O16").Value2 = 46.44
$ws2.Range("P16").Value2 = 0.016
$ws2.Range("Q16").Value2 = 46.421
$ws2.Range("R16").Value2 = 46.478
$ws2.Range("B28").Value2 = 3.991
$ws2.Range("C28").Value2 = 3.952
$ws2.Range("D28").Value2 = 3.97
$ws2.Range("E28").Value2 = 3.968
$ws2.Range("F28").Value2 = 3.975
$ws2.Range("G28").Value2 = 3.949
$ws2.Range("H28").Value2 = 3.973
$ws2.Range("I28").Value2 = 3.966
$ws2.Range("J28").Value2 = 3.948
$ws2.Range("K28").Value2 = 3.972
$ws2.Range("L28").Value2 = 3.97
$ws2.Range("M28").Value2 = 4.012
$ws2.Range("N28").Value2 = 3.959
$ws2.Range("O28").Value2 = 3.97
$ws2.Range("P28").Value2 = 0.017
$ws2.Range("Q28").Value2 = 3.948
$ws2.Range("R28").Value2 = 4.012
$ws2.Range("B29").Value2 = 0.066
$ws2.Range("C29").Value2 = 0.106
$ws2.Range("D29").Value2 = 0.08699999999999999
$ws2.Range("E29").Value2 = 0.089
$ws2.Range("F29").Value2 = 0.082
$ws2.Range("G29").Value2 = 0.109
$ws2.Range("H29").Value2 = 0.08500000000000001
$ws2.Range("I29").Value2 = 0.092
$ws2.Range("J29").Value2 = 0.109
$ws2.Range("K29").Value2 = 0.08599999999999999
$ws2.Range("L29").Value2 = 0.08699999999999999
$ws2.Range("M29").Value2 = 0.045
$ws2.Range("N29").Value2 = 0.099
$ws2.Range("O29").Value2 = 0.08799999999999999
$ws2.Range("P29").Value2 = 0.017
$ws2.Range("Q29").Value2 = 0.045
$ws2.Range("R29").Value2 = 0.109
$ws2.Range("B37").Value2 = 15.562
$ws2.Range("C37").Value2 = 15.573
$ws2.Range("D37").Value2 = 15.568
$ws2.Range("E37").Value2 = 15.568
$ws2.Range("F37").Value2 = 15.566
$ws2.Range("G37").Value2 = 15.574
$ws2.Range("H37").Value2 = 15.567
$ws2.Range("I37").Value2 = 15.569
$ws2.Range("J37").Value2 = 15.574
$ws2.Range("K37").Value2 = 15.567
$ws2.Range("L37").Value2 = 15.568
$ws2.Range("M37").Value2 = 15.556
$ws2.Range("N37").Value2 = 15.571
$ws2.Range("O37").Value2 = 15.568
$ws2.Range("Q37").Value2 = 15.556
$ws2.Range("R37").Value2 = 15.574
$ws2.Range("C38").Value2 = 5.768
$ws2.Range("E38").Value2 = 5.767
$ws2.Range("F38").Value2 = 5.766
$ws2.Range("G38").Value2 = 5.769
$ws2.Range("H38").Value2 = 5.766
$ws2.Range("I38").Value2 = 5.767
$ws2.Range("J38").Value2 = 5.769
$ws2.Range("M38").Value2 = 5.762
$ws2.Range("N38").Value2 = 5.768
$ws2.Range("Q38").Value2 = 5.762
$ws2.Range("R38").Value2 = 5.769
$ws2.Range("B39").Value2 = 4.979
$ws2.Range("C39").Value2 = 4.982
$ws2.Range("D39").Value2 = 4.981
$ws2.Range("E39").Value2 = 4.981
$ws2.Range("F39").Value2 = 4.98
$ws2.Range("G39").Value2 = 4.983
$ws2.Range("I39").Value2 = 4.981
$ws2.Range("J39").Value2 = 4.983
$ws2.Range("L39").Value2 = 4.981
$ws2.Range("M39").Value2 = 4.977
$ws2.Range("N39").Value2 = 4.982
$ws2.Range("O39").Value2 = 4.981
$ws2.Range("Q39").Value2 = 4.977
$ws2.Range("R39").Value2 = 4.983
$ws2.Range("C40").Value2 = 4.21
$ws2.Range("E40").Value2 = 4.209
$ws2.Range("I40").Value2 = 4.209
$ws2.Range("C41").Value2 = 2.065
$ws2.Range("G41").Value2 = 2.065
$ws2.Range("J41").Value2 = 2.065
$ws2.Range("R41").Value2 = 2.065
$ws2.Range("C44").Value2 = 1.087
$ws2.Range("M44").Value2 = 1.086
$ws2.Range("N44").Value2 = 1.087
$ws2.Range("Q44").Value2 = 1.086
$ws2.Range("B45").Value2 = 8.737
$ws2.Range("C45").Value2 = 8.657
$ws2.Range("D45").Value2 = 8.694000000000001
$ws2.Range("E45").Value2 = 8.69
$ws2.Range("F45").Value2 = 8.704000000000001
$ws2.Range("G45").Value2 = 8.651
$ws2.Range("H45").Value2 = 8.699999999999999
$ws2.Range("I45").Value2 = 8.685
$ws2.Range("J45").Value2 = 8.65
$ws2.Range("K45").Value2 = 8.696999999999999
$ws2.Range("L45").Value2 = 8.694000000000001
$ws2.Range("M45").Value2 = 8.779999999999999
$ws2.Range("N45").Value2 = 8.670999999999999
$ws2.Range("O45").Value2 = 8.693
$ws2.Range("P45").Value2 = 0.035
$ws2.Range("Q45").Value2 = 8.65
$ws2.Range("R45").Value2 = 8.779999999999999
$ws2.Range("B46").Value2 = 0.093
$ws2.Range("C46").Value2 = 0.149
$ws2.Range("D46").Value2 = 0.123
$ws2.Range("E46").Value2 = 0.125
$ws2.Range("F46").Value2 = 0.116
$ws2.Range("G46").Value2 = 0.153
$ws2.Range("H46").Value2 = 0.119
$ws2.Range("I46").Value2 = 0.129
$ws2.Range("J46").Value2 = 0.153
$ws2.Range("K46").Value2 = 0.121
$ws2.Range("L46").Value2 = 0.123
$ws2.Range("M46").Value2 = 0.063
$ws2.Range("N46").Value2 = 0.139
$ws2.Range("O46").Value2 = 0.124
$ws2.Range("P46").Value2 = 0.025
$ws2.Range("Q46").Value2 = 0.063
$ws2.Range("R46").Value2 = 0.153
$ws2.Range("D51").Value2 = 57.259
$ws2.Range("E51").Value2 = 57.259
$ws2.Range("H51").Value2 = 57.259
$ws2.Range("K51").Value2 = 57.259
$ws2.Range("L51").Value2 = 57.259
$ws2.Range("O51").Value2 = 57.259

$ws3 = $wb.Worksheets.Item("Edi09_3_bg_apf")
$ws3.Range("C10").Value2 = 0.44
$ws3.Range("D10").Value2 = 0.442
$ws3.Range("E10").Value2 = 0.442
$ws3.Range("F10").Value2 = 0.443
$ws3.Range("G10").Value2 = 0.439
$ws3.Range("J10").Value2 = 0.439
$ws3.Range("K10").Value2 = 0.442
$ws3.Range("L10").Value2 = 0.442
$ws3.Range("M10").Value2 = 0.448
$ws3.Range("N10").Value2 = 0.441
$ws3.Range("O10").Value2 = 0.442
$ws3.Range("Q10").Value2 = 0.439
$ws3.Range("R10").Value2 = 0.448
$ws3.Range("B11").Value2 = 0.08
$ws3.Range("C11").Value2 = 0.127
$ws3.Range("D11").Value2 = 0.105
$ws3.Range("E11").Value2 = 0.107
$ws3.Range("F11").Value2 = 0.099
$ws3.Range("G11").Value2 = 0.131
$ws3.Range("H11").Value2 = 0.102
$ws3.Range("I11").Value2 = 0.11
$ws3.Range("J11").Value2 = 0.131
$ws3.Range("K11").Value2 = 0.103
$ws3.Range("L11").Value2 = 0.105
$ws3.Range("M11").Value2 = 0.054
$ws3.Range("N11").Value2 = 0.119
$ws3.Range("O11").Value2 = 0.106
$ws3.Range("P11").Value2 = 0.021
$ws3.Range("Q11").Value2 = 0.054
$ws3.Range("R11").Value2 = 0.131
$ws3.Range("B16").Value2 = 46.447
$ws3.Range("C16").Value2 = 46.405
$ws3.Range("D16").Value2 = 46.425
$ws3.Range("E16").Value2 = 46.423
$ws3.Range("F16").Value2 = 46.43
$ws3.Range("G16").Value2 = 46.402
$ws3.Range("H16").Value2 = 46.428
$ws3.Range("I16").Value2 = 46.42
$ws3.Range("J16").Value2 = 46.402
$ws3.Range("K16").Value2 = 46.426
$ws3.Range("L16").Value2 = 46.425
$ws3.Range("M16").Value2 = 46.47
$ws3.Range("N16").Value2 = 46.413
$ws3.Range("O16").Value2 = 46.424
$ws3.Range("P16").Value2 = 0.019
$ws3.Range("Q16").Value2 = 46.402
$ws3.Range("R16").Value2 = 46.47
$ws3.Range("B28").Value2 = 3.977
$ws3.Range("C28").Value2 = 3.93
$ws3.Range("D28").Value2 = 3.952
$ws3.Range("E28").Value2 = 3.95
$ws3.Range("F28").Value2 = 3.958
$ws3.Range("G28").Value2 = 3.926
$ws3.Range("H28").Value2 = 3.955
$ws3.Range("I28").Value2 = 3.947
$ws3.Range("J28").Value2 = 3.926
$ws3.Range("K28").Value2 = 3.954
$ws3.Range("L28").Value2 = 3.952
$ws3.Range("M28").Value2 = 4.003
$ws3.Range("N28").Value2 = 3.939
$ws3.Range("O28").Value2 = 3.952
$ws3.Range("P28").Value2 = 0.021
$ws3.Range("Q28").Value2 = 3.926
$ws3.Range("R28").Value2 = 4.003
$ws3.Range("B29").Value2 = 0.08
$ws3.Range("C29").Value2 = 0.127
$ws3.Range("D29").Value2 = 0.105
$ws3.Range("E29").Value2 = 0.107
$ws3.Range("F29").Value2 = 0.099
$ws3.Range("G29").Value2 = 0.131
$ws3.Range("H29").Value2 = 0.102
$ws3.Range("I29").Value2 = 0.11
$ws3.Range("J29").Value2 = 0.131
$ws3.Range("K29").Value2 = 0.103
$ws3.Range("L29").Value2 = 0.105
$ws3.Range("M29").Value2 = 0.054
$ws3.Range("N29").Value2 = 0.119
$ws3.Range("O29").Value2 = 0.106
$ws3.Range("P29").Value2 = 0.021
$ws3.Range("Q29").Value2 = 0.054
$ws3.Range("R29").Value2 = 0.131
$ws3.Range("B37").Value2 = 15.566
$ws3.Range("C37").Value2 = 15.58
$ws3.Range("D37").Value2 = 15.573
$ws3.Range("E37").Value2 = 15.574
$ws3.Range("F37").Value2 = 15.571
$ws3.Range("G37").Value2 = 15.581
$ws3.Range("H37").Value2 = 15.572
$ws3.Range("I37").Value2 = 15.575
$ws3.Range("J37").Value2 = 15.581
$ws3.Range("K37").Value2 = 15.572
$ws3.Range("L37").Value2 = 15.573
$ws3.Range("M37").Value2 = 15.558
$ws3.Range("N37").Value2 = 15.577
$ws3.Range("O37").Value2 = 15.573
$ws3.Range("Q37").Value2 = 15.558
$ws3.Range("R37").Value2 = 15.581
$ws3.Range("B38").Value2 = 5.766
$ws3.Range("C38").Value2 = 5.771
$ws3.Range("E38").Value2 = 5.769
$ws3.Range("F38").Value2 = 5.768
$ws3.Range("G38").Value2 = 5.771
$ws3.Range("H38").Value2 = 5.768
$ws3.Range("I38").Value2 = 5.769
$ws3.Range("J38").Value2 = 5.771
$ws3.Range("M38").Value2 = 5.763
$ws3.Range("N38").Value2 = 5.77
$ws3.Range("Q38").Value2 = 5.763
$ws3.Range("R38").Value2 = 5.771
$ws3.Range("C39").Value2 = 4.984
$ws3.Range("F39").Value2 = 4.982
$ws3.Range("G39").Value2 = 4.985
$ws3.Range("H39").Value2 = 4.982
$ws3.Range("I39").Value2 = 4.983
$ws3.Range("J39").Value2 = 4.985
$ws3.Range("M39").Value2 = 4.978
$ws3.Range("N39").Value2 = 4.983
$ws3.Range("Q39").Value2 = 4.978
$ws3.Range("R39").Value2 = 4.985
$ws3.Range("C40").Value2 = 4.212
$ws3.Range("D40").Value2 = 4.21
$ws3.Range("E40").Value2 = 4.21
$ws3.Range("G40").Value2 = 4.212
$ws3.Range("H40").Value2 = 4.21
$ws3.Range("J40").Value2 = 4.212
$ws3.Range("K40").Value2 = 4.21
$ws3.Range("L40").Value2 = 4.21
$ws3.Range("M40").Value2 = 4.206
$ws3.Range("N40").Value2 = 4.211
$ws3.Range("O40").Value2 = 4.21
$ws3.Range("Q40").Value2 = 4.206
$ws3.Range("R40").Value2 = 4.212
$ws3.Range("C41").Value2 = 2.066
$ws3.Range("D41").Value2 = 2.065
$ws3.Range("E41").Value2 = 2.065
$ws3.Range("G41").Value2 = 2.066
$ws3.Range("H41").Value2 = 2.065
$ws3.Range("J41").Value2 = 2.066
$ws3.Range("K41").Value2 = 2.065
$ws3.Range("L41").Value2 = 2.065
$ws3.Range("M41").Value2 = 2.063
$ws3.Range("O41").Value2 = 2.065
$ws3.Range("Q41").Value2 = 2.063
$ws3.Range("R41").Value2 = 2.066
$ws3.Range("F44").Value2 = 1.087
$ws3.Range("H44").Value2 = 1.087
$ws3.Range("K44").Value2 = 1.087
$ws3.Range("B45").Value2 = 8.709
$ws3.Range("C45").Value2 = 8.613
$ws3.Range("D45").Value2 = 8.657999999999999
$ws3.Range("E45").Value2 = 8.654
$ws3.Range("F45").Value2 = 8.67
$ws3.Range("G45").Value2 = 8.606
$ws3.Range("H45").Value2 = 8.664
$ws3.Range("I45").Value2 = 8.647
$ws3.Range("J45").Value2 = 8.605
$ws3.Range("K45").Value2 = 8.662000000000001
$ws3.Range("L45").Value2 = 8.657999999999999
$ws3.Range("M45").Value2 = 8.760999999999999
$ws3.Range("N45").Value2 = 8.630000000000001
$ws3.Range("O45").Value2 = 8.657
$ws3.Range("P45").Value2 = 0.043
$ws3.Range("Q45").Value2 = 8.605
$ws3.Range("R45").Value2 = 8.760999999999999
$ws3.Range("B46").Value2 = 0.112
$ws3.Range("C46").Value2 = 0.18
$ws3.Range("D46").Value2 = 0.148
$ws3.Range("E46").Value2 = 0.151
$ws3.Range("F46").Value2 = 0.14
$ws3.Range("G46").Value2 = 0.184
$ws3.Range("H46").Value2 = 0.144
$ws3.Range("I46").Value2 = 0.156
$ws3.Range("J46").Value2 = 0.185
$ws3.Range("K46").Value2 = 0.145
$ws3.Range("L46").Value2 = 0.148
$ws3.Range("M46").Value2 = 0.076
$ws3.Range("N46").Value2 = 0.167
$ws3.Range("O46").Value2 = 0.149
$ws3.Range("P46").Value2 = 0.03
$ws3.Range("Q46").Value2 = 0.076
$ws3.Range("R46").Value2 = 0.185
